$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 (start time) and B2 (end time) with new raw nanosecond timer values.
$ws.Range("B1").Value = 1510873871772000000
$ws.Range("B2").Value = 1510874023621000000

# B1 switches from the "0.00" number format to the same General format as B2.
$ws.Range("B1").NumberFormat = "general"

# B3 (difference) and B4 (in seconds) hold formulas referencing B1/B2 and
# recalculate automatically.

# Move/save the active selection to B4 to match the edited view state.
$ws.Range("B4").Select() | Out-Null
